{"js": "// Updates the lattice-multiplication exercise table: each of the 15 cells\n// gets its multiplication problem (and lattice scaffolding) replaced with a\n// new problem, while keeping the existing run formatting (sz=32) intact.\nconst newCellValues = [\n  [\"50 x 33\\u000b  3    3\\u000b  ----\\u000b5|    |\\u000b0|    |\", \"34 x 45\\u000b  4    5\\u000b  ----\\u000b3|    |\\u000b4|    |\", \"37 x 28\\u000b  2    8\\u000b  ----\\u000b3|    |\\u000b7|    |\"],\n  [\"22 x 26\\u000b  2    6\\u000b  ----\\u000b2|    |\\u000b2|    |\", \"39 x 40\\u000b  4    0\\u000b  ----\\u000b3|    |\\u000b9|    |\", \"45 x 81\\u000b  8    1\\u000b  ----\\u000b4|    |\\u000b5|    |\"],\n  [\"92 x 94\\u000b  9    4\\u000b  ----\\u000b9|    |\\u000b2|    |\", \"68 x 14\\u000b  1    4\\u000b  ----\\u000b6|    |\\u000b8|    |\", \"96 x 94\\u000b  9    4\\u000b  ----\\u000b9|    |\\u000b6|    |\"],\n  [\"67 x 84\\u000b  8    4\\u000b  ----\\u000b6|    |\\u000b7|    |\", \"77 x 47\\u000b  4    7\\u000b  ----\\u000b7|    |\\u000b7|    |\", \"77 x 55\\u000b  5    5\\u000b  ----\\u000b7|    |\\u000b7|    |\"],\n  [\"43 x 99\\u000b  9    9\\u000b  ----\\u000b4|    |\\u000b3|    |\", \"39 x 98\\u000b  9    8\\u000b  ----\\u000b3|    |\\u000b9|    |\", \"80 x 79\\u000b  7    9\\u000b  ----\\u000b8|    |\\u000b0|    |\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < newCellValues.length; r++) {\n  const rowValues = newCellValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    const paragraph = cell.body.paragraphs.getFirst();\n    const range = paragraph.getRange();\n    range.insertText(rowValues[c], \"Replace\");\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Updates the lattice-multiplication exercise table: each of the 15 cells\n# gets its multiplication problem (and lattice scaffolding) replaced with a\n# new problem, while keeping the existing run formatting (sz=32) intact.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$lf = [char]11\n\n$newValues = @(\n    @((\"50 x 33\" + $lf + \"  3    3\" + $lf + \"  ----\" + $lf + \"5|    |\" + $lf + \"0|    |\"), (\"34 x 45\" + $lf + \"  4    5\" + $lf + \"  ----\" + $lf + \"3|    |\" + $lf + \"4|    |\"), (\"37 x 28\" + $lf + \"  2    8\" + $lf + \"  ----\" + $lf + \"3|    |\" + $lf + \"7|    |\")),\n    @((\"22 x 26\" + $lf + \"  2    6\" + $lf + \"  ----\" + $lf + \"2|    |\" + $lf + \"2|    |\"), (\"39 x 40\" + $lf + \"  4    0\" + $lf + \"  ----\" + $lf + \"3|    |\" + $lf + \"9|    |\"), (\"45 x 81\" + $lf + \"  8    1\" + $lf + \"  ----\" + $lf + \"4|    |\" + $lf + \"5|    |\")),\n    @((\"92 x 94\" + $lf + \"  9    4\" + $lf + \"  ----\" + $lf + \"9|    |\" + $lf + \"2|    |\"), (\"68 x 14\" + $lf + \"  1    4\" + $lf + \"  ----\" + $lf + \"6|    |\" + $lf + \"8|    |\"), (\"96 x 94\" + $lf + \"  9    4\" + $lf + \"  ----\" + $lf + \"9|    |\" + $lf + \"6|    |\")),\n    @((\"67 x 84\" + $lf + \"  8    4\" + $lf + \"  ----\" + $lf + \"6|    |\" + $lf + \"7|    |\"), (\"77 x 47\" + $lf + \"  4    7\" + $lf + \"  ----\" + $lf + \"7|    |\" + $lf + \"7|    |\"), (\"77 x 55\" + $lf + \"  5    5\" + $lf + \"  ----\" + $lf + \"7|    |\" + $lf + \"7|    |\")),\n    @((\"43 x 99\" + $lf + \"  9    9\" + $lf + \"  ----\" + $lf + \"4|    |\" + $lf + \"3|    |\"), (\"39 x 98\" + $lf + \"  9    8\" + $lf + \"  ----\" + $lf + \"3|    |\" + $lf + \"9|    |\"), (\"80 x 79\" + $lf + \"  7    9\" + $lf + \"  ----\" + $lf + \"8|    |\" + $lf + \"0|    |\")),\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Count; $c++) {\n        $cell = $t.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
